$wb = $excel.ActiveWorkbook

# The "d_faixas_idade" sheet holds the age-bracket lookup table.
$ws = $wb.Worksheets.Item("d_faixas_idade")

# Re-label the age brackets: "Menos de 20 anos" becomes "<= 20 anos" and
# "50 anos ou mais" becomes "> 50 " (trailing space kept, matching source).
$ws.Range("B2").Value = "<= 20 anos"
$ws.Range("B3").Value = "21 - 30 anos"
$ws.Range("B4").Value = "31 - 40 anos"
$ws.Range("B5").Value = "41 - 50 anos"
$ws.Range("B6").Value = "> 50 "

# A stray, empty underlined cell was left behind at B13 (new dataviz sketch
# area), which extends the used range of the sheet down to row 13.
$ws.Range("B13").Font.Underline = $true
$ws.Range("B13").Select()
